$d = $word.ActiveDocument

# Find the "SMARTREWARDS FAQ'S" heading paragraph and the final FAQ
# answer paragraph (the one ending with "...every second block 1000
# addresses will get paid.") by scanning the document paragraphs -
# using the document-level Paragraphs collection (rather than Ranges
# re-derived from Find hits) so Start/End stay reliable.
$headingIndex = -1
$lastIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*SMARTREWARDS FAQ*") {
        $headingIndex = $i
    }
    if ($t -like "*every second block 1000 addresses will get paid.*") {
        $lastIndex = $i
    }
}

if ($headingIndex -eq -1) {
    throw "Could not find the SMARTREWARDS FAQ'S heading"
}
if ($lastIndex -eq -1) {
    throw "Could not find the end of the FAQ section"
}

$startPos = $d.Paragraphs.Item($headingIndex).Range.End
$endPos = $d.Paragraphs.Item($lastIndex).Range.End

# Remove every paragraph between the heading and the trailing blank
# paragraph (inclusive of their paragraph marks), collapsing the whole
# FAQ section back down to just the heading.
$toRemove = $d.Range($startPos, $endPos)
$toRemove.Delete()
